$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param([string]$CellRef, [string]$NewValue)
    $cell = $ws.Range($CellRef)
    # Force text storage so numeric-looking strings (e.g. "1.00", "603.58")
    # are not auto-converted to numbers by Excel's General-format parser.
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    # Drop the temporary Text number-format again so the cell keeps the
    # workbook's default (unstyled) formatting, matching the source data.
    $cell.ClearFormats()
}

Set-TextCell "D2" "66.559.04"
Set-TextCell "D3" "3.234.60"
Set-TextCell "E3" "  +0.90%  "
Set-TextCell "D4" "1.00"
Set-TextCell "E4" "  -0.07%  "
Set-TextCell "D5" "603.58"
Set-TextCell "E5" "  +0.40%  "
Set-TextCell "D6" "156.18"
Set-TextCell "E6" "  -1.18%  "
Set-TextCell "E7" "  +0.03%  "
Set-TextCell "D8" "3.234.73"
Set-TextCell "E8" "  +0.87%  "
Set-TextCell "D9" "0.544"
Set-TextCell "E9" "  -1.48%  "
Set-TextCell "E10" "  +0.67%  "
Set-TextCell "D11" "5.73"
Set-TextCell "E11" "  -4.67%  "
Set-TextCell "E12" "  -2.72%  "
Set-TextCell "E13" "  +1.24%  "
Set-TextCell "D14" "38.80"
Set-TextCell "E14" "  -1.11%  "
Set-TextCell "B15" "WrappedEther"
Set-TextCell "C15" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D15" "3.304.99"
Set-TextCell "E15" "  +3.00%  "
Set-TextCell "B16" "WrappedliquidstakedEther2.0"
Set-TextCell "C16" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell "D16" "3.768.98"
Set-TextCell "E16" "  +0.97%  "
Set-TextCell "B17" "WrappedBTC"
Set-TextCell "C17" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell "D17" "66.622.73"
Set-TextCell "E17" "  -0.23%  "
Set-TextCell "D18" "7.27"
Set-TextCell "E18" "  -2.57%  "
Set-TextCell "E19" "  +1.08%  "
Set-TextCell "D20" "507.74"
Set-TextCell "E20" "  -2.10%  "
Set-TextCell "D21" "15.22"
Set-TextCell "E21" "  -1.19%  "
Set-TextCell "E22" "  -0.22%  "
Set-TextCell "D23" "8.01"
Set-TextCell "E23" "  -2.25%  "
Set-TextCell "D24" "14.59"
Set-TextCell "E24" "  -2.79%  "
Set-TextCell "D25" "86.28"
Set-TextCell "E25" "  +1.24%  "
Set-TextCell "D26" "0.166"
Set-TextCell "E26" "  +84.09%  "
Set-TextCell "D27" "1.00"
Set-TextCell "E27" "  +0.14%  "
Set-TextCell "D29" "9.02"
Set-TextCell "E29" "  -3.57%  "
Set-TextCell "E30" "  -3.18%  "
Set-TextCell "D31" "6.96"
Set-TextCell "E31" "  -1.26%  "
Set-TextCell "E32" "  -6.42%  "
Set-TextCell "E33" "  -0.17%  "
Set-TextCell "E34" "  +0.05%  "
Set-TextCell "E35" "  -6.06%  "
Set-TextCell "E36" "  -3.74%  "
Set-TextCell "D37" "0.0₃0791"
Set-TextCell "E37" "  +14.58%  "
Set-TextCell "D38" "55.32"
Set-TextCell "E38" "  +0.56%  "
Set-TextCell "D39" "492.29"
Set-TextCell "E39" "  -6.34%  "
Set-TextCell "D40" "3.16"
Set-TextCell "E40" "  +7.07%  "
Set-TextCell "D41" "0.0420"
Set-TextCell "E41" "  -1.10%  "
Set-TextCell "E42" "  +0.73%  "
Set-TextCell "E43" "  -2.46%  "
Set-TextCell "E44" "  -4.88%  "
Set-TextCell "D45" "2.943.29"
Set-TextCell "E46" "  -1.90%  "
Set-TextCell "D47" "28.13"
Set-TextCell "E47" "  -2.42%  "
Set-TextCell "E48" "  -0.92%  "
Set-TextCell "E49" "  +0.43%  "
Set-TextCell "E50" "  +0.00%  "
Set-TextCell "E51" "  -4.27%  "
